$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0299 to SCD0018
$ws.Name = "SCD0018"

# Update the TC_ID cells (B2, B3) from "DGS-314" to "SCD0018-022"
$ws.Range("B2").Value = "SCD0018-022"
$ws.Range("B3").Value = "SCD0018-022"

# Column B needs to widen to fit the new, longer TC_ID text (was bestFit for "DGS-314")
$ws.Columns.Item(2).AutoFit()

# Update the selection to match the target (activeCell B4, no topLeftCell override)
$ws.Range("B4").Select()
